# The commit removes the "H 72" record (originally row 2) from the data
# table, shifting every subsequent record up by one row and shrinking the
# used range from A1:F63 to A1:F62.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire second row (the "H 72" record); Excel shifts the rows
# below it up to fill the gap, matching the observed diff exactly.
$ws.Rows.Item(2).Delete()
